# Apply numeric cell updates to match the target workbook state.
# Values with a new value of $null are cleared (the source cell is removed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H40").Value = 1800.6666
$ws.Range("N40").ClearContents()
$ws.Range("J69").Value = 4000
$ws.Range("L69").Value = 12000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -13748
$ws.Range("H69").Value = 4000
$ws.Range("L70").Value = 78343.5
$ws.Range("K70").Value = 8178
$ws.Range("N70").Value = -78883.5
$ws.Range("M70").Value = -7908
$ws.Range("H70").Value = 9408.429
$ws.Range("I70").Value = 2726
$ws.Range("J70").Value = 26114.5
$ws.Range("H72").Value = 4000
$ws.Range("J72").Value = 4000
$ws.Range("N72").Value = -44736
$ws.Range("K72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("L72").Value = 36000
$ws.Range("I73").Value = 2726
$ws.Range("M73").Value = -7242
$ws.Range("N73").Value = -80215.5
$ws.Range("L73").Value = 78343.5
$ws.Range("H73").Value = 9408.429
$ws.Range("K73").Value = 8178
$ws.Range("J73").Value = 26114.5
$ws.Range("K76").Value = 7866.1113
$ws.Range("N76").Value = -7264.6665
$ws.Range("H76").Value = 7558.25
$ws.Range("L76").Value = 6634.6665
$ws.Range("M76").Value = -7551.1113
$ws.Range("J76").Value = 6634.6665
$ws.Range("I76").Value = 7866.1113
$ws.Range("M79").Value = -6774.1113
$ws.Range("I79").Value = 7866.1113
$ws.Range("J79").Value = 6634.6665
$ws.Range("K79").Value = 7866.1113
$ws.Range("H79").Value = 7558.25
$ws.Range("L79").Value = 6634.6665
$ws.Range("N79").Value = -8818.666499999999
$ws.Range("M86").Value = -562.75
$ws.Range("H86").Value = 6309.4546
$ws.Range("I86").Value = 1685.75
$ws.Range("K86").Value = 1685.75
$ws.Range("K89").Value = 8428.75
$ws.Range("I89").Value = 1685.75
$ws.Range("M89").Value = -2812.75
$ws.Range("H89").Value = 6309.4546
$ws.Range("I94").Value = 0
$ws.Range("H94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("K94").Value = 0
$ws.Range("K100").Value = 2365.889
$ws.Range("I100").Value = 2365.889
$ws.Range("H100").Value = 3844.7273
$ws.Range("M100").Value = -1824.889
$ws.Range("I138").Value = 7694.7334
$ws.Range("H138").Value = 10816.907
$ws.Range("M138").Value = -17944.2002
$ws.Range("N138").Value = -44444.108
$ws.Range("J138").Value = 11388.036
$ws.Range("K138").Value = 23084.2002
$ws.Range("L138").Value = 34164.108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 2185378.8
$ws.Range("K32").Value = 2185378.8
$ws.Range("N32").Value = -67577.164
$ws.Range("H32").Value = 1747094.2
$ws.Range("J32").Value = 67003.164
$ws.Range("L32").Value = 67003.164
$ws.Range("M32").Value = -2185091.8
$ws.Range("H36").Value = 3707.2942
$ws.Range("I36").Value = 6004.1665
$ws.Range("M36").Value = -5658.1665
$ws.Range("K36").Value = 6004.1665
$ws.Range("L36").Value = 2454.4546
$ws.Range("N36").Value = -3146.4546
$ws.Range("J36").Value = 2454.4546
$ws.Range("I45").Value = 8999
$ws.Range("J45").Value = 9000
$ws.Range("K45").Value = 8999
$ws.Range("N45").Value = -9754
$ws.Range("H45").Value = 8999.666999999999
$ws.Range("L45").Value = 9000
$ws.Range("M45").Value = -8622
$ws.Range("M74").Value = -1367.2173
$ws.Range("K74").Value = 2241.2173
$ws.Range("I74").Value = 2241.2173
$ws.Range("H74").Value = 13610.135
$ws.Range("M77").Value = -6838.086499999999
$ws.Range("K77").Value = 11206.0865
$ws.Range("H77").Value = 13610.135
$ws.Range("I77").Value = 2241.2173
$ws.Range("M122").Value = -8384.636200000001
$ws.Range("K122").Value = 10834.6362
$ws.Range("I122").Value = 3611.5454
$ws.Range("H122").Value = 5099.8887
$ws.Range("M132").Value = -12759.08
$ws.Range("K132").Value = 15289.08
$ws.Range("H132").Value = 11179.194
$ws.Range("I132").Value = 5096.36

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I94").Value = 3378.7896
$ws.Range("L94").Value = 3426.6667
$ws.Range("H94").Value = 3390.28
$ws.Range("M94").Value = -2927.7896
$ws.Range("N94").Value = -4328.6667
$ws.Range("K94").Value = 3378.7896
$ws.Range("J94").Value = 3426.6667
$ws.Range("L134").Value = 67231.24800000001
$ws.Range("J134").Value = 22410.416
$ws.Range("N134").Value = -72301.24800000001
$ws.Range("M134").Value = -15139.827
$ws.Range("H134").Value = 11555.2
$ws.Range("I134").Value = 5891.609
$ws.Range("K134").Value = 17674.827

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 6636.386
$ws.Range("L31").Value = 43080.645
$ws.Range("K31").Value = 6636.386
$ws.Range("N31").Value = -43670.645
$ws.Range("J31").Value = 43080.645
$ws.Range("M31").Value = -6341.386
$ws.Range("H31").Value = 15433.275
$ws.Range("L34").Value = 43080.645
$ws.Range("M34").Value = -6434.386
$ws.Range("K34").Value = 6636.386
$ws.Range("I34").Value = 6636.386
$ws.Range("N34").Value = -43484.645
$ws.Range("H34").Value = 15433.275
$ws.Range("J34").Value = 43080.645
$ws.Range("N43").Value = -7867.6665
$ws.Range("J43").Value = 7499.6665
$ws.Range("L43").Value = 7499.6665
$ws.Range("H43").Value = 7499.6665
$ws.Range("J101").Value = 7499.6665
$ws.Range("H101").Value = 7499.6665
$ws.Range("N101").Value = -13989.6665
$ws.Range("L101").Value = 7499.6665
$ws.Range("L132").Value = 33884.142
$ws.Range("J132").Value = 11294.714
$ws.Range("H132").Value = 8895.826999999999
$ws.Range("N132").Value = -38944.142
$ws.Range("M134").Value = -4154.3334
$ws.Range("H134").Value = 8702.714
$ws.Range("I134").Value = 2229.7778
$ws.Range("K134").Value = 6689.3334
$ws.Range("N141").Value = -632931.1
$ws.Range("L141").Value = 622571.1
$ws.Range("H141").Value = 564506.75
$ws.Range("J141").Value = 622571.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K4").Value = 136131312
$ws.Range("M4").Value = -136131200
$ws.Range("N4").Value = -485552.66
$ws.Range("I4").Value = 45377104
$ws.Range("L4").Value = 485328.66
$ws.Range("H4").Value = 25526472
$ws.Range("J4").Value = 161776.22
$ws.Range("J117").Value = 21296.889
$ws.Range("H117").Value = 13378.733
$ws.Range("N117").Value = -70774.667
$ws.Range("L117").Value = 63890.667
$ws.Range("L131").Value = 4465.5483
$ws.Range("N131").Value = -14545.5483
$ws.Range("H131").Value = 1468.134
$ws.Range("J131").Value = 1488.5161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I5").Value = 5499.8
$ws.Range("J5").Value = 11816.667
$ws.Range("M5").Value = -5387.8
$ws.Range("L5").Value = 11816.667
$ws.Range("H5").Value = 7868.625
$ws.Range("K5").Value = 5499.8
$ws.Range("N5").Value = -12040.667
$ws.Range("J101").Value = 18999.5
$ws.Range("H101").Value = 18999.5
$ws.Range("N101").Value = -25489.5
$ws.Range("L101").Value = 18999.5
$ws.Range("L104").Value = 36167.75
$ws.Range("N104").Value = -43155.75
$ws.Range("J104").Value = 36167.75
$ws.Range("H104").Value = 36167.75
$ws.Range("L105").Value = 173112.12
$ws.Range("H105").Value = 173112.12
$ws.Range("J105").Value = 173112.12
$ws.Range("N105").Value = -180100.12
$ws.Range("H113").Value = 102168.3
$ws.Range("K113").Value = 113960.375
$ws.Range("M113").Value = -111790.375
$ws.Range("I113").Value = 113960.375
$ws.Range("M122").Value = -7880.349999999999
$ws.Range("K122").Value = 10330.35
$ws.Range("J122").Value = 11147.5
$ws.Range("N122").Value = -38342.5
$ws.Range("I122").Value = 3443.45
$ws.Range("H122").Value = 5221.3076
$ws.Range("L122").Value = 33442.5
$ws.Range("M132").Value = -16973.7998
$ws.Range("K132").Value = 19503.7998
$ws.Range("H132").Value = 5262
$ws.Range("I132").Value = 6501.2666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J46").Value = 2999.5715
$ws.Range("L46").Value = 2999.5715
$ws.Range("H46").Value = 3079.7
$ws.Range("N46").Value = -3375.5715
$ws.Range("I61").Value = 1500
$ws.Range("H61").Value = 2928.4285
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1298
$ws.Range("K82").Value = 6603.25
$ws.Range("I82").Value = 6603.25
$ws.Range("H82").Value = 7582.6
$ws.Range("M82").Value = -6242.25
$ws.Range("H85").Value = 7582.6
$ws.Range("K85").Value = 6603.25
$ws.Range("M85").Value = -5355.25
$ws.Range("I85").Value = 6603.25
$ws.Range("H113").Value = 2928.4285
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("I113").Value = 1500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I2").Value = 1003155
$ws.Range("M2").Value = -1003043
$ws.Range("J2").Value = 901
$ws.Range("N2").Value = -1125
$ws.Range("L2").Value = 901
$ws.Range("K2").Value = 1003155
$ws.Range("H2").Value = 836112.7
$ws.Range("J15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("L62").Value = 0
$ws.Range("H62").Value = 9995
$ws.Range("J62").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H65").Value = 9995
$ws.Range("N65").ClearContents()
$ws.Range("J65").Value = 0
$ws.Range("N68").Value = -81621
$ws.Range("J68").Value = 79999
$ws.Range("H68").Value = 79999
$ws.Range("L68").Value = 79999
$ws.Range("N71").Value = -248109
$ws.Range("J71").Value = 79999
$ws.Range("L71").Value = 239997
$ws.Range("H71").Value = 79999
$ws.Range("M122").Value = -6320.5792
$ws.Range("K122").Value = 8770.5792
$ws.Range("I122").Value = 2923.5264
$ws.Range("H122").Value = 5531.2964
